# "added logout and changed application URL"
# The Login sheet's A1 cell used to hold a raw phone-number style value
# (used as a "login id"); it is now replaced with a mailto hyperlink to
# ghause000@gmail.com, mirroring the existing mailto hyperlink already
# present in B1. Column A is widened to fit the new, longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the Hyperlink look first (matches B1's style/xfId) before we
# overwrite the value, then write the new email address as the cell text.
$ws.Range("A1").Style = "Hyperlink"
$ws.Range("A1").Value = "ghause000@gmail.com"

# Turn that text into a real mailto: hyperlink, just like B1.
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:ghause000@gmail.com")

# Widen column A so the new, longer email address is fully visible.
$ws.Columns("A").AutoFit()

# Restore the sheet selection to where the author last left it.
$ws.Range("H13:I14").Select() | Out-Null
